# Daily attendance processing - 2025-10-14 09:24:15
# Re-applies the daily attendance sync: reorders "Recorded By" so the
# automated backup account is listed first, bumps a handful of session
# head-counts / coverage stats, and flips three newly-recorded sessions
# (B2D/B2E/B2F, session 15) from "Pending" to "Recorded".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) "Recorded By" cell reordering (plain text -> plain text, style s="2"
#    is preserved automatically since these values are never number- or
#    percent-like).
# ---------------------------------------------------------------------
$ws.Range("G2").Value  = "backup@backdoor.com, System, system"
$ws.Range("G4").Value  = "backup@backdoor.com, System"
$ws.Range("G5").Value  = "backup@backdoor.com, System"
$ws.Range("G29").Value = "backup@backdoor.com, System, system"
$ws.Range("G31").Value = "backup@backdoor.com, System"
$ws.Range("G32").Value = "backup@backdoor.com, System"
$ws.Range("G56").Value = "backup@backdoor.com, System, system"
$ws.Range("G58").Value = "backup@backdoor.com, System"
$ws.Range("G59").Value = "backup@backdoor.com, System"
$ws.Range("G84").Value = "backup@backdoor.com, System"
$ws.Range("G85").Value = "backup@backdoor.com, System"
$ws.Range("G110").Value = "backup@backdoor.com, System"
$ws.Range("G111").Value = "backup@backdoor.com, System"
$ws.Range("G136").Value = "backup@backdoor.com, System"
$ws.Range("G137").Value = "backup@backdoor.com, System"

# ---------------------------------------------------------------------
# 2) Students "x/y" head-count refresh (plain text fractions).
# ---------------------------------------------------------------------
$ws.Range("H9").Value   = "21/53"
$ws.Range("H36").Value  = "27/57"
$ws.Range("H37").Value  = "20/57"
$ws.Range("H63").Value  = "29/55"
$ws.Range("H87").Value  = "26/56"
$ws.Range("H90").Value  = "44/56"
$ws.Range("H92").Value  = "45/56"
$ws.Range("H116").Value = "54/55"
$ws.Range("H118").Value = "46/55"
$ws.Range("H140").Value = "51/57"
$ws.Range("H142").Value = "54/57"

# ---------------------------------------------------------------------
# 3) Overview numeric counters (K/L block).
# ---------------------------------------------------------------------
$ws.Range("L6").Value = 99
$ws.Range("L8").Value = 60

# ---------------------------------------------------------------------
# 4) Per-group statistics numeric counts (rows 18-20, columns O/Q).
# ---------------------------------------------------------------------
$ws.Range("O18").Value = 15
$ws.Range("Q18").Value = 11
$ws.Range("O19").Value = 15
$ws.Range("Q19").Value = 11
$ws.Range("O20").Value = 15
$ws.Range("Q20").Value = 11

# ---------------------------------------------------------------------
# 5) Percentage cells. These are stored as literal text (e.g. "62.3%"),
#    not as numeric percentages. Typing a "NN.N%" string straight into a
#    General-formatted cell would make Excel auto-convert it into a
#    numeric percentage, so force the cell to Text format first (exactly
#    what a user would do in the UI to paste a literal percent string),
#    then write the value.
# ---------------------------------------------------------------------
function Set-TextPercent($addr, $text) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $text
}

Set-TextPercent "L9"  "62.3%"
Set-TextPercent "L10" "71.1%"
Set-TextPercent "S15" "69.5%"
Set-TextPercent "S16" "67.8%"
Set-TextPercent "S17" "60.7%"
Set-TextPercent "R18" "57.7%"
Set-TextPercent "S18" "78.5%"
Set-TextPercent "R19" "57.7%"
Set-TextPercent "S19" "76.6%"
Set-TextPercent "R20" "57.7%"
Set-TextPercent "S20" "76.5%"

# ---------------------------------------------------------------------
# 6) Sessions that just got recorded: B2D/B2E/B2F session 15 (rows 97,
#    123, 149) flip from the "Pending" (yellow) style to the "Recorded"
#    (green) style, matching the row directly above it in each block,
#    and pick up the recorder + attendance values.
# ---------------------------------------------------------------------
function Complete-Session($row, $sourceRow, $recordedBy, $students) {
    $ws.Range("G$row").Value = $recordedBy
    $ws.Range("H$row").Value = $students
    $ws.Range("I$row").Value = "Recorded"

    # Copy the "Recorded" row's formatting (fill/font/alignment) from the
    # row above onto the whole A:I block of the newly recorded row.
    $src = $ws.Range("A$sourceRow`:I$sourceRow")
    $dst = $ws.Range("A$row`:I$row")
    $src.Copy() | Out-Null
    $dst.PasteSpecial(-4122) | Out-Null
}

Complete-Session 97  96  "System, dnasr281@gmail.com" "43/56"
Complete-Session 123 122 "System, dnasr281@gmail.com" "31/55"
Complete-Session 149 148 "System, dnasr281@gmail.com" "44/57"

$excel.CutCopyMode = $false
